# Presentacion -Grupo 2 -Proyecto 1_DHT11_ESP32_FINAL.pptx
# "se sube ppt modificado de Presentación-Grupo 2-Proyecto IOT 1 y 2"
#
# This script reproduces the textual edits from the commit:
#   1. Slide 10 title: "Conexión ESP32 + sensor DHT11(1) " ->
#      "Conexión Serial Bluetooth Terminal"
#   2-7. A handful of same-slide re-typed spans on slides 2, 3 and 6 where
#      PowerPoint merged consecutive runs that already shared identical
#      run properties (no visible text changed, only run boundaries).
#      We reproduce that by re-assigning the exact same text onto the
#      sub-range that spans the runs to be merged, which is what makes
#      PowerPoint coalesce them into a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 10 - title text replaced outright
# ---------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)
$title10 = $slide10.Shapes.Item(1)
$title10.TextFrame.TextRange.Text = "Conexión Serial Bluetooth Terminal"

# ---------------------------------------------------------------------
# 2) Slide 2 - "Introducción" body placeholder: merge re-typed runs
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$body2 = $slide2.Shapes.Item(2)
$tr2 = $body2.TextFrame.TextRange

# "onitorea las variables de temperatura " + "y humedad " + "usando "
$m1 = $tr2.Characters(16, 55)
$m1.Text = $m1.Text

# " " + "ESP32 y un sensor " + "DHT11, " + "comunicando los datos a un Smartphone o PC vía "
$m2 = $tr2.Characters(88, 73)
$m2.Text = $m2.Text

# "l " + "microcontrolador"
$m3 = $tr2.Characters(221, 18)
$m3.Text = $m3.Text

# ---------------------------------------------------------------------
# 3) Slide 3 - "Componentes Utilizados" body placeholder
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$body3 = $slide3.Shapes.Item(2)
$tr3 = $body3.TextFrame.TextRange

# "ESP32" + ")"
$m4 = $tr3.Characters(106, 6)
$m4.Text = $m4.Text

# " " + "    ("
$m5 = $tr3.Characters(113, 6)
$m5.Text = $m5.Text

# ---------------------------------------------------------------------
# 4) Slide 6 - "Video del Funcionamiento" body placeholder
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$body6 = $slide6.Shapes.Item(2)
$tr6 = $body6.TextFrame.TextRange

# "Envía datos vía Bluetooth al dispositivo " + "emparejado."
$m6 = $tr6.Characters(159, 52)
$m6.Text = $m6.Text
